# Apply updated values to range C2:M25 on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(0.02676330784007064,0.126626398853702,0.1272319464966003,1.76866632858848,1.139794440754088,1.118416158928923,0.9677211564322121,0.1582491484842379,1.444391610187438,0,0.4473300512769285),
  @(0.02514591287319234,0.124834256285439,0.1272293766838999,1.777265849831636,1.145804523967982,1.127957161975999,0.9715550296392479,0.1594774004235973,1.300306291216657,0,0.4188996723942324),
  @(0.02414555528302742,0.1237593896984919,0.1272791742192183,1.783783967672548,1.150478680072538,1.134503208182139,0.9746639311407677,0.1603161508165449,1.211928809736406,0,0.4015609944499374),
  @(0.02373609068877158,0.1233278363438401,0.1273123997652021,1.786750863919217,1.152630066833567,1.137343506281937,0.9761201776870081,0.1606792063535636,1.175938654400113,0,0.3945252112629163),
  @(0.02366799056277102,0.1232565690214074,0.127318698354685,1.787262263694181,1.153002175699982,1.137825561438916,0.976373408214819,0.1607407749369827,1.169964030734917,0,0.3933587370636431),
  @(0.02414004040819862,0.1237535433850212,0.1272795699281346,1.783822723007837,1.150506696883667,1.134540814365266,0.9746828045842051,0.1603209610583072,1.211443333321569,0,0.4014659861101464),
  @(0.02620714895996912,0.1260031996545905,0.127220405455752,1.771374143830116,1.141662023431266,1.121562998471134,0.9688861088638063,0.1586550821181412,1.394692778114461,0,0.4375029838088764),
  @(0.03020248734353714,0.1306154622780014,0.1275115484714142,1.756813449219536,1.132161452857318,1.101581915422784,0.9635311178379311,0.1560602832327547,1.754729185052156,0,0.5090972688332656),
  @(0.0331018837429724,0.1341245227246333,0.1279732271162324,1.752162328161887,1.130015487548363,1.090250896577786,0.9632948581471936,0.1545646317565357,2.019639047456849,0,0.5622573552422523),
  @(0.03441300797710056,0.1357466688723292,0.1282369823307157,1.751368436910013,1.130100281091643,1.085826732779012,0.9639975397990952,0.1539736548532176,2.140235050164563,0,0.5865622976395315),
  @(0.03490835934017156,0.136364613000481,0.1283445771743104,1.751258605394597,1.130285859568531,1.084256719874915,0.9643806843307772,0.1537627431769018,2.185913306537714,0,0.5957833499657852),
  @(0.03480172766905554,0.1362313650546128,0.1283210616740966,1.751273762520114,1.13023905203697,1.084590161308881,0.9642929524397772,0.1538075936474534,2.176075201499657,0,0.5937966662077798),
  @(0.03445378383776898,0.1357974341617592,0.1282456796336824,1.751355573079422,1.130112467850722,1.085695454369628,0.964026712209396,0.153956044789787,2.143992811715805,0,0.5873205739911498),
  @(0.03424050883599961,0.135532115898414,0.1282005105585107,1.751430552918819,1.130054944242303,1.086386203086505,0.9638788927723354,0.1540486532256082,2.124342835895334,0,0.5833560254229724),
  @(0.03301603970192701,0.134019028051128,0.1279570704728421,1.752240871715301,1.130031372937168,1.09055474826188,0.9632652860747939,0.154605054411558,2.011759451542957,0,0.5606714051207717),
  @(0.03226285053285949,0.1330973875627919,0.1278214831477875,1.753077055330948,1.13028929986082,1.093299291422156,0.9630967533290544,0.1549693025736438,1.942714535664209,0,0.5467862192500945),
  @(0.03182890047458642,0.1325697207895402,0.1277485553443753,1.753682389699208,1.130537465562753,1.094946608520843,0.9630760403779277,0.1551872221087542,1.903010011489982,0,0.538811340715327),
  @(0.03168184664968976,0.1323914817526699,0.1277247323615427,1.753908687047698,1.130638605841199,1.09551615700515,0.9630821017602713,0.1552624502857256,1.889568207071875,0,0.5361131695119212),
  @(0.03234310506370974,0.1331952459348784,0.127835393198179,1.752975163795014,1.130251506325436,1.093000014839163,0.9631068011519659,0.1549296567934633,1.950063637360472,0,0.5482631313824129),
  @(0.03455601453175916,0.1359247908709875,0.1282676118269492,1.751326359496503,1.130145476428154,1.085367942239401,0.9641017321757275,0.1539120913931242,2.15341589494443,0,0.5892222918890013),
  @(0.03599560905813348,0.137730092284329,0.1285950644073637,1.751361215770402,1.13097113206581,1.080993976261652,0.9654345168558223,0.1533221296203422,2.286382831548849,0,0.6160921784209563),
  @(0.03522788635160623,0.1367646269994793,0.128416185239562,1.751240589637561,1.130448278954759,1.083272162149328,0.96466054867944,0.1536301263806692,2.215410418557553,0,0.6017420861463307),
  @(0.03230682489138559,0.1331509973194116,0.1278290888172648,1.753020840873106,1.130268281721328,1.093135101482943,0.96310202127313,0.1549475541489045,1.946741138551204,0,0.5475953947826326),
  @(0.02912792877669546,0.129346438421031,0.1273892433794046,1.759693636847956,1.133886471045656,1.106400311938273,0.9643327389613816,0.1566901927202977,1.657259338187089,0,0.4896305062402533)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        # Columns C (3) through M (13)
        $colIndex = 3 + $j
        $ws.Cells.Item($rowIndex, $colIndex).Value = $rowVals[$j]
    }
}
